# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets
# to reflect the newly generated data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 1506
$wsExhibit.Range("F7").Value = 999
$wsExhibit.Range("F9").Value = 211
$wsExhibit.Range("F10").Value = 149
$wsExhibit.Range("F11").Value = 213
$wsExhibit.Range("F12").Value = 120
$wsExhibit.Range("F13").Value = 181
$wsExhibit.Range("F14").Value = 167

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1506
$wsAll.Range("F8").Value = 999
$wsAll.Range("F10").Value = 211
$wsAll.Range("F11").Value = 149
$wsAll.Range("F12").Value = 213
$wsAll.Range("F13").Value = 120
$wsAll.Range("F14").Value = 181
$wsAll.Range("F15").Value = 167
